# Update all planning docs
# Updated Planning Documentation to Current
#
# The document was (re)synced through a SharePoint document library, which
# stamps the package with the library's custom XML metadata parts:
#   - the content-type schema (columns/taxonomy fields for the library)
#   - the SharePoint form-template pointer
#   - the document's own property values for those columns
#
# Reproduce that via the Word object model: each SharePoint-managed custom
# XML payload is registered on the document through CustomXMLParts.Add —
# Word mints the backing itemProps part / relationship bookkeeping for us.

$d = $word.ActiveDocument

# 1) Content type schema for the document library ("Document" content type,
#    with the MediaService* / taxonomy / TaxCatchAll column definitions).
$contentTypeSchemaXml = @"
<?xml version="1.0" encoding="utf-8"?>
<ct:contentTypeSchema xmlns:ct="http://schemas.microsoft.com/office/2006/metadata/contentType" xmlns:ma="http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes" ct:_="" ma:_="" ma:contentTypeName="Document" ma:contentTypeID="0x010100305A213F757D4E4A8BAE9DDA8E12820D" ma:contentTypeVersion="10" ma:contentTypeDescription="Create a new document." ma:contentTypeScope="" ma:versionID="949ed28498f4d3f70e007772f9afb913">
  <xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:ns2="5c9051c5-e21b-4128-ad49-fdef9cda0a37" xmlns:ns3="67872523-8638-47fd-beae-94bd45b0d032" targetNamespace="http://schemas.microsoft.com/office/2006/metadata/properties" ma:root="true" ma:fieldsID="3f306e7ecddf6898cd2c1568fdea009c" ns2:_="" ns3:_="">
    <xsd:import namespace="5c9051c5-e21b-4128-ad49-fdef9cda0a37"/>
    <xsd:import namespace="67872523-8638-47fd-beae-94bd45b0d032"/>
    <xsd:element name="properties">
      <xsd:complexType>
        <xsd:sequence>
          <xsd:element name="documentManagement">
            <xsd:complexType>
              <xsd:all>
                <xsd:element ref="ns2:MediaServiceMetadata" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceFastMetadata" minOccurs="0"/>
                <xsd:element ref="ns2:lcf76f155ced4ddcb4097134ff3c332f" minOccurs="0"/>
                <xsd:element ref="ns3:TaxCatchAll" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceOCR" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceGenerationTime" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceEventHashCode" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceSearchProperties" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceObjectDetectorVersions" minOccurs="0"/>
              </xsd:all>
            </xsd:complexType>
          </xsd:element>
        </xsd:sequence>
      </xsd:complexType>
    </xsd:element>
  </xsd:schema>
  <xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" targetNamespace="5c9051c5-e21b-4128-ad49-fdef9cda0a37" elementFormDefault="qualified">
    <xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/>
    <xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
    <xsd:element name="MediaServiceMetadata" ma:index="8" nillable="true" ma:displayName="MediaServiceMetadata" ma:hidden="true" ma:internalName="MediaServiceMetadata" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceFastMetadata" ma:index="9" nillable="true" ma:displayName="MediaServiceFastMetadata" ma:hidden="true" ma:internalName="MediaServiceFastMetadata" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="lcf76f155ced4ddcb4097134ff3c332f" ma:index="11" nillable="true" ma:taxonomy="true" ma:internalName="lcf76f155ced4ddcb4097134ff3c332f" ma:taxonomyFieldName="MediaServiceImageTags" ma:displayName="Image Tags" ma:readOnly="false" ma:fieldId="{5cf76f15-5ced-4ddc-b409-7134ff3c332f}" ma:taxonomyMulti="true" ma:sspId="e764c76a-3358-464e-9b7e-b53b17345a4f" ma:termSetId="09814cd3-568e-fe90-9814-8d621ff8fb84" ma:anchorId="fba54fb3-c3e1-fe81-a776-ca4b69148c4d" ma:open="true" ma:isKeyword="false">
      <xsd:complexType>
        <xsd:sequence>
          <xsd:element ref="pc:Terms" minOccurs="0" maxOccurs="1"/>
        </xsd:sequence>
      </xsd:complexType>
    </xsd:element>
    <xsd:element name="MediaServiceOCR" ma:index="13" nillable="true" ma:displayName="Extracted Text" ma:internalName="MediaServiceOCR" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note">
          <xsd:maxLength value="255"/>
        </xsd:restriction>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceGenerationTime" ma:index="14" nillable="true" ma:displayName="MediaServiceGenerationTime" ma:hidden="true" ma:internalName="MediaServiceGenerationTime" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceEventHashCode" ma:index="15" nillable="true" ma:displayName="MediaServiceEventHashCode" ma:hidden="true" ma:internalName="MediaServiceEventHashCode" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceSearchProperties" ma:index="16" nillable="true" ma:displayName="MediaServiceSearchProperties" ma:hidden="true" ma:internalName="MediaServiceSearchProperties" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceObjectDetectorVersions" ma:index="17" nillable="true" ma:displayName="MediaServiceObjectDetectorVersions" ma:hidden="true" ma:indexed="true" ma:internalName="MediaServiceObjectDetectorVersions" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
  </xsd:schema>
  <xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" targetNamespace="67872523-8638-47fd-beae-94bd45b0d032" elementFormDefault="qualified">
    <xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/>
    <xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
    <xsd:element name="TaxCatchAll" ma:index="12" nillable="true" ma:displayName="Taxonomy Catch All Column" ma:hidden="true" ma:list="{cb7946ba-1ca9-4e3a-9109-677f16ee4972}" ma:internalName="TaxCatchAll" ma:showField="CatchAllData" ma:web="67872523-8638-47fd-beae-94bd45b0d032">
      <xsd:complexType>
        <xsd:complexContent>
          <xsd:extension base="dms:MultiChoiceLookup">
            <xsd:sequence>
              <xsd:element name="Value" type="dms:Lookup" maxOccurs="unbounded" minOccurs="0" nillable="true"/>
            </xsd:sequence>
          </xsd:extension>
        </xsd:complexContent>
      </xsd:complexType>
    </xsd:element>
  </xsd:schema>
  <xsd:schema xmlns="http://schemas.openxmlformats.org/package/2006/metadata/core-properties" xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:dc="http://purl.org/dc/elements/1.1/" xmlns:dcterms="http://purl.org/dc/terms/" xmlns:odoc="http://schemas.microsoft.com/internal/obd" targetNamespace="http://schemas.openxmlformats.org/package/2006/metadata/core-properties" elementFormDefault="qualified" attributeFormDefault="unqualified" blockDefault="#all">
    <xsd:import namespace="http://purl.org/dc/elements/1.1/" schemaLocation="http://dublincore.org/schemas/xmls/qdc/2003/04/02/dc.xsd"/>
    <xsd:import namespace="http://purl.org/dc/terms/" schemaLocation="http://dublincore.org/schemas/xmls/qdc/2003/04/02/dcterms.xsd"/>
    <xsd:element name="coreProperties" type="CT_coreProperties"/>
    <xsd:complexType name="CT_coreProperties">
      <xsd:all>
        <xsd:element ref="dc:creator" minOccurs="0" maxOccurs="1"/>
        <xsd:element ref="dcterms:created" minOccurs="0" maxOccurs="1"/>
        <xsd:element ref="dc:identifier" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="contentType" minOccurs="0" maxOccurs="1" type="xsd:string" ma:index="0" ma:displayName="Content Type"/>
        <xsd:element ref="dc:title" minOccurs="0" maxOccurs="1" ma:index="4" ma:displayName="Title"/>
        <xsd:element ref="dc:subject" minOccurs="0" maxOccurs="1"/>
        <xsd:element ref="dc:description" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="keywords" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element ref="dc:language" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="category" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element name="version" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element name="revision" minOccurs="0" maxOccurs="1" type="xsd:string">
          <xsd:annotation>
            <xsd:documentation>
                        This value indicates the number of saves or revisions. The application is responsible for updating this value after each revision.
                    </xsd:documentation>
          </xsd:annotation>
        </xsd:element>
        <xsd:element name="lastModifiedBy" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element ref="dcterms:modified" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="contentStatus" minOccurs="0" maxOccurs="1" type="xsd:string"/>
      </xsd:all>
    </xsd:complexType>
  </xsd:schema>
  <xs:schema xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" xmlns:xs="http://www.w3.org/2001/XMLSchema" targetNamespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" elementFormDefault="qualified" attributeFormDefault="unqualified">
    <xs:element name="Person">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:DisplayName" minOccurs="0"/>
          <xs:element ref="pc:AccountId" minOccurs="0"/>
          <xs:element ref="pc:AccountType" minOccurs="0"/>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="DisplayName" type="xs:string"/>
    <xs:element name="AccountId" type="xs:string"/>
    <xs:element name="AccountType" type="xs:string"/>
    <xs:element name="BDCAssociatedEntity">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:BDCEntity" minOccurs="0" maxOccurs="unbounded"/>
        </xs:sequence>
        <xs:attribute ref="pc:EntityNamespace"/>
        <xs:attribute ref="pc:EntityName"/>
        <xs:attribute ref="pc:SystemInstanceName"/>
        <xs:attribute ref="pc:AssociationName"/>
      </xs:complexType>
    </xs:element>
    <xs:attribute name="EntityNamespace" type="xs:string"/>
    <xs:attribute name="EntityName" type="xs:string"/>
    <xs:attribute name="SystemInstanceName" type="xs:string"/>
    <xs:attribute name="AssociationName" type="xs:string"/>
    <xs:element name="BDCEntity">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:EntityDisplayName" minOccurs="0"/>
          <xs:element ref="pc:EntityInstanceReference" minOccurs="0"/>
          <xs:element ref="pc:EntityId1" minOccurs="0"/>
          <xs:element ref="pc:EntityId2" minOccurs="0"/>
          <xs:element ref="pc:EntityId3" minOccurs="0"/>
          <xs:element ref="pc:EntityId4" minOccurs="0"/>
          <xs:element ref="pc:EntityId5" minOccurs="0"/>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="EntityDisplayName" type="xs:string"/>
    <xs:element name="EntityInstanceReference" type="xs:string"/>
    <xs:element name="EntityId1" type="xs:string"/>
    <xs:element name="EntityId2" type="xs:string"/>
    <xs:element name="EntityId3" type="xs:string"/>
    <xs:element name="EntityId4" type="xs:string"/>
    <xs:element name="EntityId5" type="xs:string"/>
    <xs:element name="Terms">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:TermInfo" minOccurs="0" maxOccurs="unbounded"/>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="TermInfo">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:TermName" minOccurs="0"/>
          <xs:element ref="pc:TermId" minOccurs="0"/>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="TermName" type="xs:string"/>
    <xs:element name="TermId" type="xs:string"/>
  </xs:schema>
</ct:contentTypeSchema>
"@

# 2) SharePoint "use the standard document-library forms" pointer.
$formTemplatesXml = @"
<?xml version="1.0" encoding="utf-8"?>
<?mso-contentType ?>
<FormTemplates xmlns="http://schemas.microsoft.com/sharepoint/v3/contenttype/forms">
  <Display>DocumentLibraryForm</Display>
  <Edit>DocumentLibraryForm</Edit>
  <New>DocumentLibraryForm</New>
</FormTemplates>
"@

# 3) The document's own values for the library's managed-metadata columns.
$documentPropertiesXml = @"
<?xml version="1.0" encoding="utf-8"?>
<p:properties xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls">
  <documentManagement>
    <TaxCatchAll xmlns="67872523-8638-47fd-beae-94bd45b0d032" xsi:nil="true"/>
    <lcf76f155ced4ddcb4097134ff3c332f xmlns="5c9051c5-e21b-4128-ad49-fdef9cda0a37">
      <Terms xmlns="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
    </lcf76f155ced4ddcb4097134ff3c332f>
  </documentManagement>
</p:properties>
"@

$d.CustomXMLParts.Add($contentTypeSchemaXml) | Out-Null
$d.CustomXMLParts.Add($formTemplatesXml) | Out-Null
$d.CustomXMLParts.Add($documentPropertiesXml) | Out-Null

Write-Output ("CustomXMLParts.Count = " + $d.CustomXMLParts.Count)
